# Applies the "Added final changes to cambridge maths challenge" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected (no password) - unprotect so we can edit, then
# re-protect at the end to restore the original state.
$ws.Unprotect()

# 1. Add a "Graduate" button (I/J columns) to several NPC rows that didn't have it yet.
$graduateRows = @(5, 9, 13, 17, 21, 25, 32)
foreach ($r in $graduateRows) {
    $ws.Cells.Item($r, 9).Value = "Graduate"
    $ws.Cells.Item($r, 10).Value = "scriptevent graduation:junior"
}

# 2. Rename "Guild Master" to "Guild Leader" throughout column C (scaleNpc rows).
$used = $ws.UsedRange
foreach ($row in $used.Rows) {
    $cell = $ws.Cells.Item($row.Row, 3)
    $current = $cell.Value()
    if ($current -eq "Guild Master") {
        $cell.Value = "Guild Leader"
    }
}

# 3. Replace the "§a" placeholder marker with "§1" in a few dialogue strings (column D).
$dCells = @(59, 70, 82)
foreach ($r in $dCells) {
    $cell = $ws.Cells.Item($r, 4)
    $current = $cell.Value()
    if ($current -ne $null) {
        $cell.Value = $current.Replace("§a", "§1")
    }
}

# 4. Add "Repeat Chat" buttons (K/L columns) so players can re-trigger certain dialogues.
$ws.Cells.Item(65, 11).Value = "Repeat Chat"
$ws.Cells.Item(65, 12).Value = "dialogue open @e[tag=ratioNpc] @p ratioNpc8"

$ws.Cells.Item(79, 11).Value = "Repeat Chat"
$ws.Cells.Item(79, 12).Value = "dialogue open @e[tag=fractionNpc] @p fractionNpc8"

# 5. Update the saved view/selection to match the author's last position.
$ws.Application.ActiveWindow.ScrollRow = 54
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("L65").Select()

# Restore sheet protection to its original (no password) state.
$ws.Protect($null, $true, $true, $true)
